$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.82741333333333
$ws.Range("N2").Value = 95.48223999999999
$ws.Range("O2").Value = 0.114390792932228
$ws.Range("P2").Value = 0.114390792932228
$ws.Range("Q2").Value = 635.1409433221688
$ws.Range("R2").Value = 5716.268489899519
$ws.Range("S2").Value = 0.001342643854372281
$ws.Range("T2").Value = 0.001342643854372281
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 85.46317833333335
$ws.Range("N3").Value = 256.389535
$ws.Range("O3").Value = 0.307162904935779
$ws.Range("P3").Value = 0.307162904935779
$ws.Range("Q3").Value = 1705.484612822576
$ws.Range("R3").Value = 15349.36151540318
$ws.Range("S3").Value = 0.003605276054406734
$ws.Range("T3").Value = 0.003605276054406734
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2478306666667
$ws.Range("N4").Value = 366.743492
$ws.Range("O4").Value = 0.4393704929064738
$ws.Range("P4").Value = 0.4393704929064738
$ws.Range("Q4").Value = 2439.551140255469
$ws.Range("R4").Value = 21955.96026229922
$ws.Range("S4").Value = 0.005157041724878153
$ws.Range("T4").Value = 0.005157041724878153
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.95578266666667
$ws.Range("H5").Value = 59.867348
$ws.Range("I5").Value = 0.0117373419656925
$ws.Range("J5").Value = 0.0117373419656925
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 38.69562533333333
$ws.Range("N5").Value = 116.086876
$ws.Range("O5").Value = 0.1390758092255191
$ws.Range("P5").Value = 0.1390758092255191
$ws.Range("Q5").Value = 772.2014893027609
$ws.Range("R5").Value = 6949.813403724847
$ws.Range("S5").Value = 0.001632380332035329
$ws.Range("T5").Value = 0.001632380332035329
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.82741333333333
$ws.Range("N6").Value = 95.48223999999999
$ws.Range("O6").Value = 0.114390792932228
$ws.Range("P6").Value = 0.114390792932228
$ws.Range("Q6").Value = 52112.4033568519
$ws.Range("R6").Value = 469011.6302116672
$ws.Range("S6").Value = 0.1101620023701666
$ws.Range("T6").Value = 0.1101620023701666
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 85.46317833333335
$ws.Range("N7").Value = 256.389535
$ws.Range("O7").Value = 0.307162904935779
$ws.Range("P7").Value = 0.307162904935779
$ws.Range("Q7").Value = 139932.5661441929
$ws.Range("R7").Value = 1259393.095297736
$ws.Range("S7").Value = 0.29580772887561
$ws.Range("T7").Value = 0.29580772887561
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1637.343343333333
$ws.Range("H8").Value = 4912.03003
$ws.Range("I8").Value = 0.9630320723052701
$ws.Range("J8").Value = 0.9630320723052702
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.2478306666667
$ws.Range("N8").Value = 366.743492
$ws.Range("O8").Value = 0.4393704929064738
$ws.Range("P8").Value = 0.4393704929064738
$ws.Range("Q8").Value = 200161.6717790072
$ws.Range("R8").Value = 1801455.046011065
$ws.Range("S8").Value = 0.4231278762935095
$ws.Range("T8").Value = 0.4231278762935095
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1637.343343333333
$ws.Range("H9").Value = 4912.03003
$ws.Range("I9").Value = 0.9630320723052701
$ws.Range("J9").Value = 0.9630320723052702
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.69562533333333
$ws.Range("N9").Value = 116.086876
$ws.Range("O9").Value = 0.1390758092255191
$ws.Range("P9").Value = 0.1390758092255191
$ws.Range("Q9").Value = 63358.02455565402
$ws.Range("R9").Value = 570222.2210008862
$ws.Range("S9").Value = 0.1339344647659841
$ws.Range("T9").Value = 0.1339344647659841
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.82741333333333
$ws.Range("N10").Value = 95.48223999999999
$ws.Range("O10").Value = 0.114390792932228
$ws.Range("P10").Value = 0.114390792932228
$ws.Range("Q10").Value = 557.0058105939911
$ws.Range("R10").Value = 5013.052295345919
$ws.Range("S10").Value = 0.001177471608950154
$ws.Range("T10").Value = 0.001177471608950154
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 17.50081933333334
$ws.Range("H11").Value = 52.502458
$ws.Range("I11").Value = 0.01029341242216722
$ws.Range("J11").Value = 0.01029341242216722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 85.46317833333335
$ws.Range("N11").Value = 256.389535
$ws.Range("O11").Value = 0.307162904935779
$ws.Range("P11").Value = 0.307162904935779
$ws.Range("Q11").Value = 1495.675643664115
$ws.Range("R11").Value = 13461.08079297703
$ws.Range("S11").Value = 0.003161754461294916
$ws.Range("T11").Value = 0.003161754461294916
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 17.50081933333334
$ws.Range("H12").Value = 52.502458
$ws.Range("I12").Value = 0.01029341242216722
$ws.Range("J12").Value = 0.01029341242216722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 122.2478306666667
$ws.Range("N12").Value = 366.743492
$ws.Range("O12").Value = 0.4393704929064738
$ws.Range("P12").Value = 0.4393704929064738
$ws.Range("Q12").Value = 2139.43719838926
$ws.Range("R12").Value = 19254.93478550334
$ws.Range("S12").Value = 0.004522621689617231
$ws.Range("T12").Value = 0.004522621689617232
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 17.50081933333334
$ws.Range("H13").Value = 52.502458
$ws.Range("I13").Value = 0.01029341242216722
$ws.Range("J13").Value = 0.01029341242216722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 38.69562533333333
$ws.Range("N13").Value = 116.086876
$ws.Range("O13").Value = 0.1390758092255191
$ws.Range("P13").Value = 0.1390758092255191
$ws.Range("Q13").Value = 677.2051479490232
$ws.Range("R13").Value = 6094.846331541208
$ws.Range("S13").Value = 0.001431564662304917
$ws.Range("T13").Value = 0.001431564662304917
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 25.39612333333333
$ws.Range("H14").Value = 76.18836999999999
$ws.Range("I14").Value = 0.01493717330687017
$ws.Range("J14").Value = 0.01493717330687017
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.82741333333333
$ws.Range("N14").Value = 95.48223999999999
$ws.Range("O14").Value = 0.114390792932228
$ws.Range("P14").Value = 0.114390792932228
$ws.Range("Q14").Value = 808.292914394311
$ws.Range("R14").Value = 7274.636229548799
$ws.Range("S14").Value = 0.001708675098738989
$ws.Range("T14").Value = 0.001708675098738989
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 25.39612333333333
$ws.Range("H15").Value = 76.18836999999999
$ws.Range("I15").Value = 0.01493717330687017
$ws.Range("J15").Value = 0.01493717330687017
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 85.46317833333335
$ws.Range("N15").Value = 256.389535
$ws.Range("O15").Value = 0.307162904935779
$ws.Range("P15").Value = 0.307162904935779
$ws.Range("Q15").Value = 2170.433417411994
$ws.Range("R15").Value = 19533.90075670795
$ws.Range("S15").Value = 0.004588145544467417
$ws.Range("T15").Value = 0.004588145544467417
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 25.39612333333333
$ws.Range("H16").Value = 76.18836999999999
$ws.Range("I16").Value = 0.01493717330687017
$ws.Range("J16").Value = 0.01493717330687017
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 122.2478306666667
$ws.Range("N16").Value = 366.743492
$ws.Range("O16").Value = 0.4393704929064738
$ws.Range("P16").Value = 0.4393704929064738
$ws.Range("Q16").Value = 3104.620984843115
$ws.Range("R16").Value = 27941.58886358804
$ws.Range("S16").Value = 0.006562953198468969
$ws.Range("T16").Value = 0.006562953198468969
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 25.39612333333333
$ws.Range("H17").Value = 76.18836999999999
$ws.Range("I17").Value = 0.01493717330687017
$ws.Range("J17").Value = 0.01493717330687017
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 38.69562533333333
$ws.Range("N17").Value = 116.086876
$ws.Range("O17").Value = 0.1390758092255191
$ws.Range("P17").Value = 0.1390758092255191
$ws.Range("Q17").Value = 982.718873425791
$ws.Range("R17").Value = 8844.469860832118
$ws.Range("S17").Value = 0.002077399465194792
$ws.Range("T17").Value = 0.002077399465194792
